$p = $ppt.ActivePresentation

# --- Slide 1: Title slide ---
# Title "MailScanner Installation" -> "Mail Gateway – Part 2"
# The old title text moves down into the Subtitle placeholder.
$s1 = $p.Slides.Item(1)
$title = $s1.Shapes.Item(1)
$subtitle = $s1.Shapes.Item(2)

$title.TextFrame.TextRange.Text = "Mail Gateway " + [char]0x2013 + " Part 2"
$subtitle.TextFrame.TextRange.Text = "MailScanner Installation"

# --- Slide 2: "After you run" slide ---
$s2 = $p.Slides.Item(2)
$s2Title = $s2.Shapes.Item(1)
$s2Title.TextFrame.TextRange.Text = "After you run ./install.sh"

Write-Output "Done."
